# Update countries & provincias Spain
# - refresh the "Datos actualizados ..." timestamp
# - refresh COVID numbers for a handful of countries (rows 14, 22, 71, 124)
# - two pairs of adjacent countries swapped rank (new totals re-order them):
#     row 35/36: Peru/Pakistan -> Pakistan/Peru
#     row 45/46: Panama/Finlandia -> Finlandia/Panama
#     row 59/60: Croacia/Marruecos -> Marruecos/Croacia

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 11:52"

# row -> Country, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$rows = @{
    14  = @("Suiza",      23404, 124, 9800, 12695, 391,  14, 909)
    22  = @("Israel",      9755, 351,  864,  8812, 165,   6,  79)
    35  = @("Pakistan",    4414, 151,  572,  3779,  31,   2,  63)
    36  = @("Peru",        4342,   0, 1333,  2888, 113,   0, 121)
    45  = @("Finlandia",   2605, 118,  300,  2265,  82,   0,  40)
    46  = @("Panama",      2528,   0,   16,  2449, 101,   0,  63)
    59  = @("Marruecos",   1346,  71,  103,  1147,   1,   3,  96)
    60  = @("Croacia",     1343,   0,  179,  1145,  36,   0,  19)
    71  = @("Kuwait",       910,  55,  111,   798,  21,   0,   1)
    124 = @("Brunei",       135,   0,   92,    42,   3,   0,   1)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
    $ws.Cells.Item($r, 8).Value = $vals[7]
}
